$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the six "rohit..." e-mail values to their "abc..." replacements.
$ws.Range("D46").Value = "abc123@gmail.com"
$ws.Range("D47").Value = "abc123@yahoo.com"
$ws.Range("D48").Value = "abc123@rediff.com"
$ws.Range("D50").Value = "abc@gmail.com"
$ws.Range("D53").Value = "abc123@@gmail.com"
$ws.Range("D54").Value = "abc123@gmail"

# Update the ZIP code literal value.
$ws.Range("D37").Value = 110010

# Move the selection / active cell to F54 and clear the frozen top-left cell.
$ws.Range("F54").Select()
